$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new D value, new E value). $null means "leave unchanged".
$updates = @{
    2  = @("64.044.65", "  -0.79%  ")
    3  = @("3.402.67",  "  -1.10%  ")
    4  = @($null,       "  -0.02%  ")
    5  = @("570.69",    "  -0.39%  ")
    6  = @("162.70",    "  +2.47%  ")
    7  = @($null,       "  -0.01%  ")
    8  = @("3.402.07",  "  -1.08%  ")
    9  = @("0.547",     "  -4.58%  ")
    10 = @($null,       "  +1.47%  ")
    11 = @($null,       "  -1.80%  ")
    12 = @("0.420",     "  -4.07%  ")
    13 = @("3.992.44",  "  -1.03%  ")
    14 = @($null,       "  +0.60%  ")
    15 = @("26.85",     "  -2.33%  ")
    16 = @($null,       "  -1.26%  ")
    17 = @("64.055.46", "  -0.93%  ")
    18 = @("3.393.27",  "  -1.17%  ")
    19 = @("6.11",      "  -0.69%  ")
    20 = @("13.46",     "  -1.76%  ")
    21 = @("372.52",    "  -1.33%  ")
    22 = @($null,       "  -1.44%  ")
    23 = @("0.998",     "  -0.16%  ")
    24 = @("70.13",     $null)
    25 = @($null,       "  -3.26%  ")
    26 = @($null,       "  -4.09%  ")
    27 = @("9.49",      "  -4.17%  ")
    28 = @($null,       "  -0.65%  ")
    29 = @("1.00",      "  -0.06%  ")
    30 = @($null,       "  +1.03%  ")
    31 = @("1.38",      "  -3.52%  ")
    34 = @("22.75",     "  -1.88%  ")
    35 = @("7.02",      "  +0.69%  ")
    36 = @($null,       "  -5.47%  ")
    37 = @("159.58",    "  -0.71%  ")
    38 = @("0.855",     "  +8.14%  ")
    39 = @($null,       "  -3.05%  ")
    40 = @("0.0724",    "  -2.73%  ")
    41 = @($null,       "  -1.52%  ")
    42 = @("42.77",     "  -0.33%  ")
    43 = @("6.46",      "  +0.54%  ")
    44 = @("2.724.00",  "  -5.57%  ")
    45 = @("25.88",     "  +0.59%  ")
    46 = @($null,       "  -3.50%  ")
    47 = @($null,       "  -1.65%  ")
    48 = @("2.40",      "  +0.90%  ")
    49 = @("327.61",    "  +2.32%  ")
    50 = @($null,       "  -3.10%  ")
    51 = @($null,       "  -1.80%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    if ($null -ne $dVal) {
        $cell = $ws.Range("D$row")
        # Price column holds plain numeric-looking text (e.g. "570.69") in
        # the source data. Force text formatting first so Excel's COM layer
        # doesn't auto-convert the literal into a real number, then restore
        # the default "Normal" style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.Style = "Normal"
    }
    if ($null -ne $eVal) {
        $ws.Range("E$row").Value = $eVal
    }
}
